$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PCB version bump 1.1 -> 1.2
$ws.Range("A4").Value = "PCB version: 1.2"

# Row 15: R4 resistor changed from 1 ohm to 7.5 ohm
$ws.Range("E15").Value = "PAC100007508FA1000"
$ws.Range("F15").Value = "Res 7.5 ohm 1W 1% Axial high temperature"

# Row 23: U7,U8,U9 splitter part changed from QCN-27D+ to QCN-27+
$ws.Range("E23").Value = "QCN-27+"
$ws.Range("F23").Value = "Power Splitter/Combiner, 2 Way-90°, 50Ω, 1700 to 2700 Mhz, LTCC "
$ws.Range("I23").Value = "139-QCN-27"

# Update the saved selection to match the author's final cursor position
$null = $ws.Range("I23").Select()
